$wb = $excel.ActiveWorkbook

# --- Set up sheets -----------------------------------------------------
# The original workbook has a single sheet ("Sheet1") that becomes
# "FaceNet512" (keeps sheetId=1). Two new sheets "ArcFace" (sheetId=2)
# and "VGGFace" (sheetId=3) are inserted *before* it, in that order.
$orig = $wb.Worksheets.Item(1)

$arcface = $wb.Worksheets.Add($orig)
$arcface.Name = "ArcFace"

$vggface = $wb.Worksheets.Add($null, $arcface)
$vggface.Name = "VGGFace"

$facenet = $wb.Worksheets.Item("Sheet1")
$facenet.Name = "FaceNet512"

# --- Helper data ---------------------------------------------------------
# Columns: A Metric | B Value (Weighted) | C Value (Micro) | D Value(Macro)

function Fill-Sheet($ws, $acc, $prec, $rec, $f1) {
    $ws.Range("A1").Value = "Metric"
    $ws.Range("B1").Value = "Value (Weighted)"
    $ws.Range("C1").Value = "Value (Micro)"
    $ws.Range("D1").Value = "Value(Macro)"

    $ws.Range("A2").Value = "Accuracy"
    $ws.Range("B2").Value = $acc[0]
    $ws.Range("C2").Value = $acc[1]
    $ws.Range("D2").Value = $acc[2]

    $ws.Range("A3").Value = "Precision"
    $ws.Range("B3").Value = $prec[0]
    $ws.Range("C3").Value = $prec[1]
    $ws.Range("D3").Value = $prec[2]

    $ws.Range("A4").Value = "Recall"
    $ws.Range("B4").Value = $rec[0]
    $ws.Range("C4").Value = $rec[1]
    $ws.Range("D4").Value = $rec[2]

    $ws.Range("A5").Value = "F1-Score"
    $ws.Range("B5").Value = $f1[0]
    $ws.Range("C5").Value = $f1[1]
    $ws.Range("D5").Value = $f1[2]

    $ws.Columns.Item(1).ColumnWidth = 10.08984375
    $ws.Columns.Item(2).ColumnWidth = 15.26953125
    $ws.Columns.Item(3).ColumnWidth = 12
    $ws.Columns.Item(4).ColumnWidth = 12.08984375
}

# --- ArcFace --------------------------------------------------------------
$ws = $wb.Worksheets.Item("ArcFace")
Fill-Sheet $ws `
    @(0.995305, 0.995305, 0.995305) `
    @(0.99688600000000005, 0.995305, 0.66331700000000005) `
    @(0.995305, 0.995305, 0.66347100000000003) `
    @(0.99609000000000003, 0.995305, 0.66339000000000004)
[void]$ws.Range("E7").Select()

# --- VGGFace ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("VGGFace")
Fill-Sheet $ws `
    @(0.99843499999999996, 0.99843499999999996, 0.99843499999999996) `
    @(1, 0.99843499999999996, 0.66666700000000001) `
    @(0.99843499999999996, 0.99843499999999996, 0.66498299999999999) `
    @(0.99921599999999999, 0.99843499999999996, 0.66582300000000005)
[void]$ws.Range("E3").Select()

# --- FaceNet512 ---------------------------------------------------------
$ws = $wb.Worksheets.Item("FaceNet512")
Fill-Sheet $ws `
    @(0.995305, 0.995305, 0.995305) `
    @(0.99843499999999996, 0.995305, 0.66498299999999999) `
    @(0.995305, 0.995305, 0.66347100000000003) `
    @(0.99686699999999995, 0.995305, 0.66422599999999998)
[void]$ws.Range("E3").Select()

# --- Active sheet / tab selection ------------------------------------------
$ws = $wb.Worksheets.Item("FaceNet512")
[void]$ws.Select()
